$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# shares_outstanding (H) and fixed_ticker (I) are constant across all data rows
$ws.Range("H2:H44").Value = 2850792605
$ws.Range("I2:I44").Value = "ORCL"

# open/close/high/low price series (D/E/F/G) per row
$ws.Range("D2").Value = 37.24306400146951
$ws.Range("E2").Value = 36.6819953918457
$ws.Range("F2").Value = 38.20367935666982
$ws.Range("G2").Value = 35.0752958233462
$ws.Range("D3").Value = 37.39920613925572
$ws.Range("E3").Value = 34.37928771972656
$ws.Range("F3").Value = 38.59352515110287
$ws.Range("G3").Value = 34.21720291631113
$ws.Range("D4").Value = 31.10816509049254
$ws.Range("E4").Value = 30.92834663391113
$ws.Range("F4").Value = 32.98338764468962
$ws.Range("G4").Value = 30.17483442750364
$ws.Range("D5").Value = 33.60449193892197
$ws.Range("E5").Value = 31.40373611450196
$ws.Range("F5").Value = 33.86239221787135
$ws.Range("G5").Value = 31.03407784556851
$ws.Range("D6").Value = 32.06840351101668
$ws.Range("E6").Value = 35.3140869140625
$ws.Range("F6").Value = 36.25499036721334
$ws.Range("G6").Value = 31.97344943795278
$ws.Range("D7").Value = 34.70927109987086
$ws.Range("E7").Value = 35.46306610107422
$ws.Range("F7").Value = 35.66234508736126
$ws.Range("G7").Value = 32.99373593928244
$ws.Range("D8").Value = 35.63698320959823
$ws.Range("E8").Value = 34.15863037109375
$ws.Range("F8").Value = 36.01961434301524
$ws.Range("G8").Value = 33.66294784837945
$ws.Range("D9").Value = 35.01555101248438
$ws.Range("E9").Value = 33.56639099121094
$ws.Range("F9").Value = 35.91472763303917
$ws.Range("G9").Value = 32.85926947600123
$ws.Range("D10").Value = 37.46639449288043
$ws.Range("E10").Value = 39.09651184082031
$ws.Range("F10").Value = 41.18236117208305
$ws.Range("G10").Value = 37.02819104851925
$ws.Range("D11").Value = 40.0659259440413
$ws.Range("E11").Value = 44.13236999511719
$ws.Range("F11").Value = 45.63748194577347
$ws.Range("G11").Value = 38.49920145909834
$ws.Range("D12").Value = 44.48437112272137
$ws.Range("E12").Value = 42.71736526489258
$ws.Range("F12").Value = 46.94934510326362
$ws.Range("G12").Value = 42.09007899423445
$ws.Range("D13").Value = 43.50675167648554
$ws.Range("E13").Value = 41.93678283691406
$ws.Range("F13").Value = 44.97915131497105
$ws.Range("G13").Value = 41.68842733130614
$ws.Range("D14").Value = 45.37690538274636
$ws.Range("E14").Value = 40.73770523071289
$ws.Range("F14").Value = 47.62081873421386
$ws.Range("G14").Value = 39.63355609709607
$ws.Range("D15").Value = 42.09558244994729
$ws.Range("E15").Value = 39.39531326293945
$ws.Range("F15").Value = 43.51724386315139
$ws.Range("G15").Value = 38.06305994800678
$ws.Range("D16").Value = 43.40966200385802
$ws.Range("E16").Value = 46.28209686279297
$ws.Range("F16").Value = 46.77579580162358
$ws.Range("G16").Value = 42.48509834981842
$ws.Range("D17").Value = 44.55779645444244
$ws.Range("E17").Value = 40.69143676757812
$ws.Range("F17").Value = 45.26077375586883
$ws.Range("G17").Value = 38.21299939594923
$ws.Range("D18").Value = 47.50294181276178
$ws.Range("E18").Value = 48.59777069091797
$ws.Range("F18").Value = 49.20400050507813
$ws.Range("G18").Value = 46.14571490382573
$ws.Range("D19").Value = 45.95247105063323
$ws.Range("E19").Value = 51.77832794189453
$ws.Range("F19").Value = 52.23276296763118
$ws.Range("G19").Value = 45.34352631327278
$ws.Range("D20").Value = 47.39780185961966
$ws.Range("E20").Value = 50.21757888793945
$ws.Range("F20").Value = 51.36739277270608
$ws.Range("G20").Value = 47.04190824612959
$ws.Range("D21").Value = 51.54150928394218
$ws.Range("E21").Value = 48.56249618530273
$ws.Range("F21").Value = 51.9356559121076
$ws.Range("G21").Value = 48.16834955713731
$ws.Range("D22").Value = 45.78635957449623
$ws.Range("E22").Value = 44.49738311767578
$ws.Range("F22").Value = 47.72903497238984
$ws.Range("G22").Value = 36.56095542973775
$ws.Range("D23").Value = 49.28789106867275
$ws.Range("E23").Value = 51.12878799438477
$ws.Range("F23").Value = 51.77633983639333
$ws.Range("G23").Value = 47.09546892407997
$ws.Range("D24").Value = 53.29580810725422
$ws.Range("E24").Value = 55.46034240722656
$ws.Range("F24").Value = 58.15439373029037
$ws.Range("G24").Value = 50.97334823713411
$ws.Range("D25").Value = 54.28579316144116
$ws.Range("E25").Value = 60.33931350708008
$ws.Range("F25").Value = 61.74775420121157
$ws.Range("G25").Value = 54.10856922502986
$ws.Range("D26").Value = 60.81365251966963
$ws.Range("E26").Value = 65.70121765136719
$ws.Range("F26").Value = 68.93150848175571
$ws.Range("G26").Value = 60.64511554466227
$ws.Range("D27").Value = 74.8343238858025
$ws.Range("E27").Value = 73.19808197021484
$ws.Range("F27").Value = 79.95931543414696
$ws.Range("G27").Value = 71.37377690916358
$ws.Range("D28").Value = 84.24668533804817
$ws.Range("E28").Value = 82.23631286621094
$ws.Range("F28").Value = 87.08762253295926
$ws.Range("G28").Value = 79.97111540070151
$ws.Range("D29").Value = 86.84127566709687
$ws.Range("E29").Value = 82.59818267822266
$ws.Range("F29").Value = 100.716552410588
$ws.Range("G29").Value = 82.56977028654167
$ws.Range("D30").Value = 72.2293673121553
$ws.Range("E30").Value = 78.64616394042969
$ws.Range("F30").Value = 80.33829581677409
$ws.Range("G30").Value = 69.50103783846555
$ws.Range("D31").Value = 69.68703043605703
$ws.Range("E31").Value = 66.68081665039062
$ws.Range("F31").Value = 71.42395703092936
$ws.Range("G31").Value = 60.84970073820917
$ws.Range("D32").Value = 71.0531900320821
$ws.Range("E32").Value = 58.54314804077149
$ws.Range("F32").Value = 76.12431019388195
$ws.Range("G32").Value = 58.26514629499589
$ws.Range("D33").Value = 80.02295385674852
$ws.Range("E33").Value = 78.76069641113281
$ws.Range("F33").Value = 82.46073794023771
$ws.Range("G33").Value = 75.21482664842642
$ws.Range("D34").Value = 84.04659973625427
$ws.Range("E34").Value = 89.86892700195312
$ws.Range("F34").Value = 89.94630193583984
$ws.Range("G34").Value = 79.34617949639355
$ws.Range("D35").Value = 101.8609898649105
$ws.Range("E35").Value = 115.6619491577148
$ws.Range("F35").Value = 123.8687173993402
$ws.Range("G35").Value = 101.7638698818279
$ws.Range("D36").Value = 117.9391190746662
$ws.Range("E36").Value = 103.2320556640625
$ws.Range("F36").Value = 124.1864479898365
$ws.Range("G36").Value = 100.649303494892
$ws.Range("D37").Value = 113.5384514091064
$ws.Range("E37").Value = 103.1305160522461
$ws.Range("F37").Value = 114.7905327357973
$ws.Range("G37").Value = 97.09508881746588
$ws.Range("D38").Value = 109.6675026262743
$ws.Range("E38").Value = 123.3464813232422
$ws.Range("F38").Value = 130.3774601759717
$ws.Range("G38").Value = 108.3811115060655
$ws.Range("D39").Value = 115.7845283323148
$ws.Range("E39").Value = 139.1029968261719
$ws.Range("F39").Value = 143.161819791043
$ws.Range("G39").Value = 115.5973476696015
$ws.Range("D40").Value = 138.7547871346323
$ws.Range("E40").Value = 168.3432922363281
$ws.Range("F40").Value = 171.8899729464154
$ws.Range("G40").Value = 136.2849617492717
$ws.Range("D41").Value = 183.4651931344328
$ws.Range("E41").Value = 164.9988708496094
$ws.Range("F41").Value = 196.3569718867232
$ws.Range("G41").Value = 162.9294575209077
$ws.Range("D42").Value = 165.7121281340302
$ws.Range("E42").Value = 138.7733612060547
$ws.Range("F42").Value = 165.8808658331603
$ws.Range("G42").Value = 135.7459727318952
$ws.Range("D43").Value = 163.5380963361446
$ws.Range("E43").Value = 217.7884826660156
$ws.Range("F43").Value = 227.3415664356535
$ws.Range("G43").Value = 161.9741472802459
$ws.Range("D44").Value = 221.6154127447244
$ws.Range("E44").Value = 280.7527770996094
$ws.Range("F44").Value = 345.1210845253727
$ws.Range("G44").Value = 218.4109669663144
